$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.445.65"
$ws.Range("E2").Value = "  -2.14%  "
Set-TextValue $ws.Range("D3") "3.692.01"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws.Range("D5") "686.86"
$ws.Range("E5").Value = "  -1.76%  "
Set-TextValue $ws.Range("D6") "160.14"
$ws.Range("E6").Value = "  -5.54%  "
Set-TextValue $ws.Range("D7") "3.689.90"
$ws.Range("E7").Value = "  -2.86%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -5.82%  "
$ws.Range("E10").Value = "  -8.43%  "
Set-TextValue $ws.Range("D11") "7.19"
$ws.Range("E11").Value = "  -4.38%  "
Set-TextValue $ws.Range("D12") "0.437"
$ws.Range("E12").Value = "  -9.06%  "
$ws.Range("E13").Value = "  -6.58%  "
Set-TextValue $ws.Range("D14") "4.314.94"
$ws.Range("E14").Value = "  -2.86%  "
Set-TextValue $ws.Range("D15") "32.52"
$ws.Range("E15").Value = "  -10.19%  "
Set-TextValue $ws.Range("D16") "3.690.27"
$ws.Range("E16").Value = "  -3.18%  "
Set-TextValue $ws.Range("D17") "69.491.44"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("E18").Value = "  -0.95%  "
Set-TextValue $ws.Range("D19") "15.97"
$ws.Range("E19").Value = "  -9.24%  "
Set-TextValue $ws.Range("D20") "6.47"
$ws.Range("E20").Value = "  -10.33%  "
Set-TextValue $ws.Range("D21") "470.23"
$ws.Range("E21").Value = "  -8.20%  "
$ws.Range("E22").Value = "  -4.35%  "
$ws.Range("E23").Value = "  -9.20%  "
Set-TextValue $ws.Range("D24") "79.56"
$ws.Range("E24").Value = "  -4.76%  "
Set-TextValue $ws.Range("D25") "3.838.38"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("E26").Value = "  +0.07%  "
Set-TextValue $ws.Range("D27") "0.0000125"
$ws.Range("E27").Value = "  -11.07%  "
Set-TextValue $ws.Range("D28") "11.01"
$ws.Range("E28").Value = "  -12.91%  "
$ws.Range("E29").Value = "  -10.10%  "
$ws.Range("E30").Value = "  -8.99%  "
Set-TextValue $ws.Range("D31") "1.75"
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D32") "6.66"
$ws.Range("E32").Value = "  -8.71%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D33") "2.02"
$ws.Range("E33").Value = "  -10.63%  "
Set-TextValue $ws.Range("D34") "0.999"
$ws.Range("E34").Value = "  -0.09%  "
Set-TextValue $ws.Range("D35") "26.81"
$ws.Range("E35").Value = "  -7.83%  "
Set-TextValue $ws.Range("D36") "0.161"
$ws.Range("E36").Value = "  -6.23%  "
Set-TextValue $ws.Range("D37") "8.22"
$ws.Range("E37").Value = "  -11.70%  "
Set-TextValue $ws.Range("D38") "6.17"
$ws.Range("E38").Value = "  -7.30%  "
Set-TextValue $ws.Range("D39") "2.27"
$ws.Range("E39").Value = "  -3.41%  "
Set-TextValue $ws.Range("D41") "0.0905"
$ws.Range("E41").Value = "  -10.03%  "
$ws.Range("E42").Value = "  -0.05%  "
Set-TextValue $ws.Range("D43") "167.65"
$ws.Range("E43").Value = "  +2.17%  "
$ws.Range("E44").Value = "  -6.69%  "
Set-TextValue $ws.Range("D45") "47.94"
$ws.Range("E45").Value = "  -2.75%  "
$ws.Range("E46").Value = "  -13.70%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D47") "1.31"
$ws.Range("E47").Value = "  -3.61%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D48") "28.74"
$ws.Range("E48").Value = "  -3.73%  "
$ws.Range("E49").Value = "  -3.53%  "
$ws.Range("E50").Value = "  -9.07%  "
Set-TextValue $ws.Range("D51") "373.97"
$ws.Range("E51").Value = "  -11.36%  "
